# Make reference url format the default for the rendered QR code
#
# This removes the two "app link" columns/fields that are no longer used
# (modules__school__playStoreLink, modules__school__appStoreLink) from the
# "Connectors" table and from the "Default Values" sheet, and switches the
# active tab/selection over to the "Default Values" sheet.

$wb = $excel.ActiveWorkbook

$wsConnectors = $wb.Worksheets.Item("Connectors")
$wsDefaults   = $wb.Worksheets.Item("Default Values")

# --- Connectors sheet: shrink Table1 from 14 to 12 columns -----------------
# Resizing the table drops the two trailing ListColumns
# (modules__school__playStoreLink / modules__school__appStoreLink), which in
# turn updates the table ref, the sheet dimension and the row spans.
$tbl = $wsConnectors.ListObjects.Item(1)
$tbl.Resize($wsConnectors.Range("A1:L2"))

# Clear (value + formatting) the now-orphaned M1/N1 header cells so they
# disappear from the sheet entirely instead of lingering as empty cells.
$wsConnectors.Range("M1:N1").Clear()

# --- Default Values sheet: blank out the matching header cells -------------
# These two cells keep their header style but lose their text value.
$wsDefaults.Range("I1:J1").ClearContents()

# --- View state: "Default Values" becomes the active / selected tab --------
$wsConnectors.Activate()
$wsConnectors.Range("L4").Select()

$wsDefaults.Activate()
$wsDefaults.Range("D4").Select()
